# Weekly data refresh: the data block in columns D:T (Fecha ... Kg/unidad)
# for rows 2-14 is re-shuffled to new row positions (new prices / dates
# pulled in, old rows reordered). Columns A,B,C stay the same for every
# row already. The mapping below (new row -> source row whose D:T values
# it should receive) was derived from the target OOXML diff:
#   new row 2  <- old row 7
#   new row 3  <- old row 8
#   new row 4  <- old row 9
#   new row 5  <- old row 6
#   new row 6  <- old row 13
#   new row 7  <- old row 14
#   new row 8  <- old row 11
#   new row 9  <- old row 12
#   new row 10 <- old row 4
#   new row 11 <- old row 2
#   new row 12 <- old row 3
#   new row 13 <- old row 10
#   new row 14 <- old row 5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 14
$firstCol = 4   # column D (Fecha)
$lastCol = 20   # column T (Kg / unidad)

# mapping: new row number -> source (old) row number
$mapping = @{
    2  = 7
    3  = 8
    4  = 9
    5  = 6
    6  = 13
    7  = 14
    8  = 11
    9  = 12
    10 = 4
    11 = 2
    12 = 3
    13 = 10
    14 = 5
}

# 1) Snapshot all the current (old) values for D:T across rows 2-14
#    before any writes happen, since this is a permutation and source
#    rows would otherwise get clobbered before being read.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the snapshot values into their new row positions.
for ($newRow = $firstRow; $newRow -le $lastRow; $newRow++) {
    $srcRow = $mapping[$newRow]
    $rowVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value = $rowVals[$c]
    }
}
